$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.438.17'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.870.40'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.87'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7059'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3156'
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07885'
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.67'
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08016'
$ws.Range("E11").Value = '  -3.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.892.88'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.13'
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7042'
$ws.Range("E15").Value = '  -1.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.495'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.513.14'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008365'
$ws.Range("E18").Value = '  -3.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '256.20'
$ws.Range("E19").Value = '  +5.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.140.88'
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.625'
$ws.Range("E23").Value = '  -2.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1556'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.073'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.14'
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.80'
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.501'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.337'
$ws.Range("E30").Value = '  -2.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.256'
$ws.Range("E31").Value = '  -2.17%  '
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05324'
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.899'
$ws.Range("E34").Value = '  -2.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7487'
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.716'
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01877'
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.264.14'
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.750'
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8986'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.06'
$ws.Range("E42").Value = '  -3.46%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.003'
$ws.Range("E43").Value = '  -7.95%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.81'
$ws.Range("E44").Value = '  -3.58%  '
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000128'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.038.67'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.811'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5196'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.519'
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4327'
$ws.Range("E51").Value = '  -1.20%  '
